{"js": "// Replace the date line and each \"a\u00f7b=c, d\" answer cell with its new value.\n// Every \"old\" string below is unique in the document, so body.search()\n// (matchCase, exact substring) unambiguously finds the single run to update.\nconst replacements = [\n  [\"2025-04-26 Saturday\", \"2025-04-27 Sunday\"],\n  [\"80\u00f73=26, 2\", \"39\u00f79=4, 3\"],\n  [\"15\u00f78=1, 7\", \"12\u00f79=1, 3\"],\n  [\"94\u00f72=47, 0\", \"39\u00f72=19, 1\"],\n  [\"53\u00f74=13, 1\", \"95\u00f74=23, 3\"],\n  [\"36\u00f78=4, 4\", \"49\u00f74=12, 1\"],\n  [\"69\u00f72=34, 1\", \"13\u00f72=6, 1\"],\n  [\"97\u00f75=19, 2\", \"82\u00f72=41, 0\"],\n  [\"34\u00f72=17, 0\", \"18\u00f73=6, 0\"],\n  [\"84\u00f76=14, 0\", \"11\u00f77=1, 4\"],\n  [\"75\u00f77=10, 5\", \"44\u00f74=11, 0\"],\n  [\"15\u00f79=1, 6\", \"11\u00f77=1, 4\"],\n  [\"85\u00f74=21, 1\", \"24\u00f78=3, 0\"],\n  [\"85\u00f77=12, 1\", \"27\u00f79=3, 0\"],\n  [\"63\u00f77=9, 0\", \"98\u00f78=12, 2\"],\n  [\"59\u00f73=19, 2\", \"90\u00f79=10, 0\"],\n  [\"71\u00f73=23, 2\", \"15\u00f75=3, 0\"],\n  [\"12\u00f72=6, 0\", \"79\u00f72=39, 1\"],\n  [\"36\u00f75=7, 1\", \"13\u00f77=1, 6\"],\n  [\"29\u00f75=5, 4\", \"42\u00f75=8, 2\"],\n  [\"81\u00f72=40, 1\", \"48\u00f74=12, 0\"],\n  [\"43\u00f78=5, 3\", \"68\u00f78=8, 4\"],\n  [\"10\u00f77=1, 3\", \"55\u00f79=6, 1\"],\n  [\"50\u00f79=5, 5\", \"67\u00f79=7, 4\"],\n  [\"28\u00f78=3, 4\", \"26\u00f72=13, 0\"],\n  [\"46\u00f77=6, 4\", \"56\u00f72=28, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"a\u00f7b=c, d\" answer cell with its new value.\n# Every \"old\" string is unique in the document, so Find/Replace (MatchCase,\n# whole document, ReplaceAll) unambiguously rewrites exactly one run each.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"2025-04-26 Saturday\", $false, $true, $false, $false, $false, $true, 1, $false, \"2025-04-27 Sunday\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"80\u00f73=26, 2\", $false, $true, $false, $false, $false, $true, 1, $false, \"39\u00f79=4, 3\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"15\u00f78=1, 7\", $false, $true, $false, $false, $false, $true, 1, $false, \"12\u00f79=1, 3\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"94\u00f72=47, 0\", $false, $true, $false, $false, $false, $true, 1, $false, \"39\u00f72=19, 1\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"53\u00f74=13, 1\", $false, $true, $false, $false, $false, $true, 1, $false, \"95\u00f74=23, 3\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"36\u00f78=4, 4\", $false, $true, $false, $false, $false, $true, 1, $false, \"49\u00f74=12, 1\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"69\u00f72=34, 1\", $false, $true, $false, $false, $false, $true, 1, $false, \"13\u00f72=6, 1\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"97\u00f75=19, 2\", $false, $true, $false, $false, $false, $true, 1, $false, \"82\u00f72=41, 0\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"34\u00f72=17, 0\", $false, $true, $false, $false, $false, $true, 1, $false, \"18\u00f73=6, 0\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"84\u00f76=14, 0\", $false, $true, $false, $false, $false, $true, 1, $false, \"11\u00f77=1, 4\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"75\u00f77=10, 5\", $false, $true, $false, $false, $false, $true, 1, $false, \"44\u00f74=11, 0\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"15\u00f79=1, 6\", $false, $true, $false, $false, $false, $true, 1, $false, \"11\u00f77=1, 4\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"85\u00f74=21, 1\", $false, $true, $false, $false, $false, $true, 1, $false, \"24\u00f78=3, 0\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"85\u00f77=12, 1\", $false, $true, $false, $false, $false, $true, 1, $false, \"27\u00f79=3, 0\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"63\u00f77=9, 0\", $false, $true, $false, $false, $false, $true, 1, $false, \"98\u00f78=12, 2\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"59\u00f73=19, 2\", $false, $true, $false, $false, $false, $true, 1, $false, \"90\u00f79=10, 0\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"71\u00f73=23, 2\", $false, $true, $false, $false, $false, $true, 1, $false, \"15\u00f75=3, 0\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"12\u00f72=6, 0\", $false, $true, $false, $false, $false, $true, 1, $false, \"79\u00f72=39, 1\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"36\u00f75=7, 1\", $false, $true, $false, $false, $false, $true, 1, $false, \"13\u00f77=1, 6\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"29\u00f75=5, 4\", $false, $true, $false, $false, $false, $true, 1, $false, \"42\u00f75=8, 2\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"81\u00f72=40, 1\", $false, $true, $false, $false, $false, $true, 1, $false, \"48\u00f74=12, 0\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"43\u00f78=5, 3\", $false, $true, $false, $false, $false, $true, 1, $false, \"68\u00f78=8, 4\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"10\u00f77=1, 3\", $false, $true, $false, $false, $false, $true, 1, $false, \"55\u00f79=6, 1\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"50\u00f79=5, 5\", $false, $true, $false, $false, $false, $true, 1, $false, \"67\u00f79=7, 4\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"28\u00f78=3, 4\", $false, $true, $false, $false, $false, $true, 1, $false, \"26\u00f72=13, 0\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"46\u00f77=6, 4\", $false, $true, $false, $false, $false, $true, 1, $false, \"56\u00f72=28, 0\", 2)\n"}
